$d = $word.ActiveDocument

$pairs = @(
    @("674÷6=112, 2", "277÷7=39, 4"),
    @("946÷9=105, 1", "968÷5=193, 3"),
    @("915÷6=152, 3", "966÷4=241, 2"),
    @("870÷5=174, 0", "883÷7=126, 1"),
    @("179÷4=44, 3", "346÷9=38, 4"),
    @("637÷4=159, 1", "381÷6=63, 3"),
    @("859÷7=122, 5", "183÷5=36, 3"),
    @("489÷8=61, 1", "166÷8=20, 6"),
    @("314÷5=62, 4", "298÷3=99, 1"),
    @("163÷6=27, 1", "885÷8=110, 5"),
    @("628÷2=314, 0", "883÷9=98, 1"),
    @("303÷2=151, 1", "522÷5=104, 2"),
    @("749÷7=107, 0", "567÷6=94, 3"),
    @("986÷8=123, 2", "216÷6=36, 0"),
    @("857÷5=171, 2", "951÷7=135, 6"),
    @("232÷9=25, 7", "514÷8=64, 2"),
    @("862÷4=215, 2", "216÷7=30, 6"),
    @("435÷4=108, 3", "360÷9=40, 0"),
    @("886÷7=126, 4", "639÷4=159, 3"),
    @("903÷5=180, 3", "367÷5=73, 2"),
    @("397÷9=44, 1", "277÷2=138, 1"),
    @("763÷4=190, 3", "608÷5=121, 3"),
    @("304÷4=76, 0", "631÷4=157, 3"),
    @("140÷6=23, 2", "133÷2=66, 1"),
    @("631÷2=315, 1", "244÷9=27, 1")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
